# "uml: Tweak SD chapter"
#
# Re-balances the sequence-diagram callback slide: the method-call labels
# (write()/getText()/getAuthor()) and the lifeline boxes (:Chapter/:Book)
# shrink from 24pt to 20pt and get nudged/resized to line up with the
# arrows, and the two lifeline header boxes move to sit right above their
# lifelines (and, since PowerPoint draws z-order == document order, they
# get bumped to the end of the shape stack so they paint on top).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 34 "Rectangle 33" (activation box under :Chapter) ---------------
$shp = $s.Shapes.Item(4)
$shp.TextFrame.TextRange.Font.Size = 20

# --- Shape 37 "Line 16" (dashed red return line, bottom-left) --------------
$shp = $s.Shapes.Item(7)
$shp.TextFrame.TextRange.Font.Size = 20

# --- Shape 38 "Line 9" (blue lifeline under :Book) --------------------------
$shp = $s.Shapes.Item(8)
$shp.TextFrame.TextRange.Font.Size = 20

# --- Shape 39 "Rectangle 38" (activation box on :Book lifeline) ------------
$shp = $s.Shapes.Item(9)
$shp.Width = 20.2501968503937
$shp.Height = 122.7194094488189
$shp.TextFrame.TextRange.Font.Size = 20

# --- Shape 40 "Line 16" (dashed blue return line) ---------------------------
$shp = $s.Shapes.Item(10)
$shp.Left = 270.28043307086614
$shp.Top = 349.6298031496063
$shp.TextFrame.TextRange.Font.Size = 20

# --- Shape 41 "Rectangle 40" (small activation box, self-call) -------------
$shp = $s.Shapes.Item(11)
$shp.TextFrame.TextRange.Font.Size = 20

# --- Shape 42 "Line 16" (solid blue self-call return line) -----------------
$shp = $s.Shapes.Item(12)
$shp.TextFrame.TextRange.Font.Size = 20

# --- Shape "TextBox 1" (write() label) --------------------------------------
$shp = $s.Shapes.Item(14)
$shp.Height = 31.504763779527558
$shp.TextFrame.TextRange.Font.Size = 20

# --- Shape "TextBox 15" (getText() label) -----------------------------------
$shp = $s.Shapes.Item(15)
$shp.Left = 371.3095669291339
$shp.Top = 197.59255905511813
$shp.Width = 89.43405511811024
$shp.Height = 31.504763779527558
$shp.TextFrame.TextRange.Font.Size = 20

# --- Shape "TextBox 16" (getAuthor() label) ---------------------------------
$shp = $s.Shapes.Item(16)
$shp.Left = 300.320905511811
$shp.Top = 244.1090157480315
$shp.Width = 111.42759842519685
$shp.Height = 31.504763779527558
$shp.TextFrame.TextRange.Font.Size = 20

# --- Shapes 31 (:Chapter) and 32 (:Book) lifeline header boxes -------------
# Resize/reposition them to hug their lifelines, then send both to the
# front (end of z-order) like the authored deck does.
$chapter = $s.Shapes.Item(1)
$chapter.Left = 468.2141338582677
$chapter.Top = 152.69578740157482
$chapter.Width = 101.34145669291338
$chapter.Height = 29.959094488188974
$chapter.TextFrame.TextRange.Font.Size = 20
$chapter.TextFrame.BottomInset = 0.39370078740157477
$chapter.TextFrame.TopInset = 0.39370078740157477

$book = $s.Shapes.Item(2)
$book.Left = 204.96610236220474
$book.Top = 154.99216535433072
$book.Width = 100.46153543307086
$book.Height = 29.959094488188974
$book.TextFrame.TextRange.Font.Size = 20
$book.TextFrame.BottomInset = 0.39370078740157477
$book.TextFrame.TopInset = 0.39370078740157477

$s.Shapes.Item(1).ZOrder(0)
$s.Shapes.Item(1).ZOrder(0)
